$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F5").Value = "a"
$ws.Range("F6").Select()
